$p = $ppt.ActivePresentation

# Add a new slide at the end (position 7) using the "Title and Content" layout
# (ppLayoutText = 2), matching the layout used by the deck's other content slides.
$s = $p.Slides.Add(7, 2)

# --- Title placeholder ---
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Ciclo de vida de un script en Unity"
$title.LanguageID = "es-CO"

# --- Content placeholder: numbered list ---
$content = $s.Shapes.Item(2).TextFrame

# Match the indent/margin PowerPoint uses for a numbered (arabic period) list.
$lvl1 = $content.Ruler.Levels.Item(1)
$lvl1.LeftMargin = 40.5
$lvl1.FirstMargin = -40.5

$tr = $content.TextRange

# Paragraph 1: Inicialización
$tr.Text = "Inicialización"
$tr.LanguageID = "es-ES"
$para1 = $tr.Paragraphs(1, 1)
$para1.ParagraphFormat.Bullet.Font.Name = "+mj-lt"
$para1.ParagraphFormat.Bullet.Type = 2

# Paragraph 2: Activación
$r2 = $tr.InsertAfter("`rActivación")
$r2.LanguageID = "es-ES"
$para2 = $tr.Paragraphs(2, 1)
$para2.ParagraphFormat.Bullet.Font.Name = "+mj-lt"
$para2.ParagraphFormat.Bullet.Type = 2

# Paragraph 3: Actualización por Frame (two runs - "Frame" marked as a different word)
$r3 = $tr.InsertAfter("`rActualización por ")
$r3.LanguageID = "es-ES"
$para3 = $tr.Paragraphs(3, 1)
$para3.ParagraphFormat.Bullet.Font.Name = "+mj-lt"
$para3.ParagraphFormat.Bullet.Type = 2
$r3b = $para3.InsertAfter("Frame")
$r3b.LanguageID = "es-ES"

# Paragraph 4: Desactivación
$r4 = $tr.InsertAfter("`rDesactivación")
$r4.LanguageID = "es-ES"
$para4 = $tr.Paragraphs(4, 1)
$para4.ParagraphFormat.Bullet.Font.Name = "+mj-lt"
$para4.ParagraphFormat.Bullet.Type = 2

# Paragraph 5: Finalización
$r5 = $tr.InsertAfter("`rFinalización")
$r5.LanguageID = "es-ES"
$para5 = $tr.Paragraphs(5, 1)
$para5.ParagraphFormat.Bullet.Font.Name = "+mj-lt"
$para5.ParagraphFormat.Bullet.Type = 2

Write-Output "Added slide with $($s.Shapes.Count) shapes; deck now has $($p.Slides.Count) slides"
